# Ajuste Modulo Gestion Cotizaciones, Solicitud Modificar, Vista Adendas
#
# The "Servicio" sheet's "Unidad" column (D) and the whole "Tiempo Ejecucion"
# block (años / meses / dias, columns E:G) are removed. The trailing
# "Cantidad" column (old H) shifts left into the new column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Servicio")

# Remove columns D (Unidad) through G (Dias) entirely - this drops the
# header cells, the merged "Tiempo Ejecucion" cell, the validation lists
# tied to those columns, and shifts the old "Cantidad" column (H) left
# so it becomes the new column D.
$ws.Columns("D:G").Delete()

# Reflect the cursor position left behind on the Servicio sheet without
# disturbing which sheet/tab is active in the workbook.
$ws.Range("D14").Select()
$wb.Worksheets.Item("Bien").Activate()
